# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.867.20"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "'3.317.67"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'558.94"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'185.78"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'3.313.63"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "'0.576"
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("D10").Value = "'0.176"
$ws.Range("E10").Value = "  -6.03%  "
$ws.Range("E11").Value = "  -1.70%  "
$ws.Range("D12").Value = "'45.87"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "'3.850.24"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "'8.41"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "'573.36"
$ws.Range("E16").Value = "  -9.51%  "
$ws.Range("D17").Value = "'65.868.67"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("E18").Value = "  +0.54%  "
$ws.Range("D19").Value = "'3.314.70"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "'10.83"
$ws.Range("E21").Value = "  -4.59%  "
$ws.Range("D22").Value = "'0.890"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "'4.98"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "'98.27"
$ws.Range("E25").Value = "  -8.54%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").Value = "'9.34"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("D30").Value = "'30.48"
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  +6.57%  "
$ws.Range("D32").Value = "'3.69"
$ws.Range("E32").Value = "  -9.49%  "
$ws.Range("D33").Value = "'557.29"
$ws.Range("E33").Value = "  +5.38%  "
$ws.Range("D34").Value = "'10.82"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "'3.742.19"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "'55.59"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("D39").Value = "'33.70"
$ws.Range("E39").Value = "  +2.32%  "
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("D41").Value = "'0.0₃0684"
$ws.Range("E41").Value = "  -5.96%  "
$ws.Range("D42").Value = "'3.12"
$ws.Range("E42").Value = "  -7.03%  "
$ws.Range("D43").Value = "'2.58"
$ws.Range("E43").Value = "  -6.13%  "
$ws.Range("D44").Value = "'3.33"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  -2.08%  "
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").Value = "'2.97"
$ws.Range("E47").Value = "  -13.12%  "
$ws.Range("D48").Value = "'0.127"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").Value = "'2.51"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").Value = "'124.88"
$ws.Range("E51").Value = "  +2.91%  "
